$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305865526199341
$ws.Range("B1").Value = 3.975365400314331
$ws.Range("C1").Value = 3.752973556518555
$ws.Range("D1").Value = 3.101295232772827
$ws.Range("E1").Value = 1.041487574577332
